# Geekbench workbook update — "results for Raspberry Pi 3"
#
# Adds Geekbench/benchmark results to the "2020" sheet:
#   - Row 4  (Raspberry Pi 1B 1.2): fills in several previously-empty metrics
#   - Row 5  (Raspberry Pi 3): fills in the JetStream/Speedometer/Octane/etc.
#     results, and removes the grey placeholder shading now that real data
#     is present
#   - Row 34/35: two new benchmarked machines appended to the table
#     (MacBook Air 2020 M1, iPad air 4)

$xlRight = -4152

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2020")

# ---------------------------------------------------------------------
# Row 4 - Raspberry Pi 1B 1.2
# ---------------------------------------------------------------------

# Label column gets right-aligned (was center-aligned)
$ws.Range("C4").HorizontalAlignment = $xlRight

$ws.Range("M4").Value = 255

$ws.Range("N4").Value = 155
$ws.Range("N4").HorizontalAlignment = $xlRight

$ws.Range("O4").Value = " - "

$ws.Range("P4").Value = 1.7
$ws.Range("P4").HorizontalAlignment = $xlRight

$ws.Range("Q4").Value = 0.678
$ws.Range("Q4").HorizontalAlignment = $xlRight

$ws.Range("T4").Value = 3631
$ws.Range("T4").HorizontalAlignment = $xlRight

# ---------------------------------------------------------------------
# Row 5 - Raspberry Pi 3
# ---------------------------------------------------------------------

# The row was shaded grey as a "results pending" placeholder; clear the
# fill now that the real numbers are being filled in (copy the
# (unformatted) format of a plain cell over so the grey highlight goes
# away instead of just blanking the values).
$blank = $ws.Range("R4")
$placeholderRange = $ws.Range("O5:U5")
$blank.Copy() | Out-Null
$placeholderRange.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("N5").Value = 1586
$ws.Range("P5").Value = 9.9
$ws.Range("Q5").Value = 4.9
$ws.Range("R5").Value = 8.421
$ws.Range("S5").Value = 1.86
$ws.Range("U5").Value = 22988

# ---------------------------------------------------------------------
# New rows 34 & 35 - two additional benchmarked machines
# ---------------------------------------------------------------------

$ws.Range("A34").Value = "MacBook Air 2020 M1"
$ws.Range("B34").Value = "M1"
$ws.Range("C34").Value = "8 x 3200"
$ws.Range("D34").Value = "7 core apple"
$ws.Range("O34").Value = 238
$ws.Range("Q34").Value = 227
$ws.Range("S34").Value = 1690

$ws.Range("A35").Value = "iPad air 4"
$ws.Range("B35").Value = "A14"
$ws.Range("N35").Value = 57779
$ws.Range("S35").Value = 1188
$ws.Range("U35").Value = 473

# ---------------------------------------------------------------------
# Restore the active selection to O5 (the cell that now carries the new
# Raspberry Pi 3 data)
# ---------------------------------------------------------------------

$ws.Range("O5").Select() | Out-Null
